$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "your wildest dre" + _GoBack bookmark + "ams with " -> one run
#    "your wildest dreams with "
# ---------------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("your wildest dre") | Out-Null
$dreStart = $r1.Start

$r2 = $d.Content
$r2.Find.Execute("ams with ") | Out-Null
$amsEnd = $r2.End

# Pin the run boundaries on either side of the phrase so the merge caused by
# the text edit below cannot spill into the neighbouring runs.
$d.Bookmarks.Add("ZZ_DREAM_BEFORE", $d.Range($dreStart, $dreStart)) | Out-Null
$d.Bookmarks.Add("ZZ_DREAM_AFTER", $d.Range($amsEnd, $amsEnd)) | Out-Null

# Re-write the phrase (this also removes the _GoBack bookmark, which sits
# inside the replaced range) and merge it into a single run.
$rng = $d.Content
$rng.Find.Execute("your wildest dreams with ") | Out-Null
$rng.Text = "your wildest dreams with X"
$rng2 = $d.Content
$rng2.Find.Execute("your wildest dreams with X") | Out-Null
$rng2.Text = "your wildest dreams with "

$d.Bookmarks("ZZ_DREAM_BEFORE").Delete()
$d.Bookmarks("ZZ_DREAM_AFTER").Delete()

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) "hotels &lodges" -> "hotels &" / " " / "lodges" (3 runs, space inserted)
# ---------------------------------------------------------------------------
$rLux = $d.Content
$rLux.Find.Execute("Luxurious ") | Out-Null
$luxEnd = $rLux.End

$d.Bookmarks.Add("ZZ_HOTELS_BEFORE", $d.Range($luxEnd, $luxEnd)) | Out-Null

$rHotels = $d.Content
$rHotels.Find.Execute("hotels &lodges") | Out-Null
$rHotels.Text = "hotels & lodges"

$rAmp = $d.Content
$rAmp.Find.Execute("hotels &") | Out-Null
$ampEnd = $rAmp.End

$d.Bookmarks.Add("ZZ_HOTELS_SPLIT1", $d.Range($ampEnd, $ampEnd)) | Out-Null
$d.Bookmarks.Add("ZZ_HOTELS_SPLIT2", $d.Range($ampEnd + 1, $ampEnd + 1)) | Out-Null

$d.Bookmarks("ZZ_HOTELS_BEFORE").Delete()
$d.Bookmarks("ZZ_HOTELS_SPLIT1").Delete()
$d.Bookmarks("ZZ_HOTELS_SPLIT2").Delete()

# ---------------------------------------------------------------------------
# 3 & 4) Remove the stale proofErr spell-check markers wrapping "Chobe" and
#         "Boma" by rebuilding those two list paragraphs from scratch.
# ---------------------------------------------------------------------------
function Rebuild-Paragraph($searchText, $newText) {
    $found = $d.Content
    $found.Find.Execute($searchText) | Out-Null
    $paraStart = $found.Paragraphs(1).Range.Start
    $paraIndex = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Start -eq $paraStart) {
            $paraIndex = $i
            break
        }
    }
    $d.Paragraphs($paraIndex).Range.Delete()
    $d.Paragraphs($paraIndex).Range.InsertBefore($newText + "`r")
}

Rebuild-Paragraph "Chobe day trip" "Chobe day trip"
Rebuild-Paragraph "Boma dinner and drum show" "Boma dinner and drum show"

Write-Output "done"
